# Updates Suggested Lists for eMedications.03, eProcedures.03 and eSituation.11
# Applies the "Intentional self-harm / X83" removal (NEMSIS XSD regex rule)
# to the eSituation.11/12 Provider Impression suggested-list workbook:
#   - "Recommendation" sheet: note appended (in red) to the existing entry,
#     the entry itself struck through, and the "Most Recent Update" banner
#     bumped to 11/20/2020.
#   - "Change Log" sheet: trailing periods added to two existing notes and
#     a new row logging the 11/20/2020 removal.

$wb = $excel.ActiveWorkbook
$wsRec = $wb.Worksheets.Item("Recommendation")
$wsLog = $wb.Worksheets.Item("Change Log")

# ---------------------------------------------------------------------------
# Recommendation sheet
# ---------------------------------------------------------------------------

# Banner note at the top of the sheet: bump the "most recent update" date.
$wsRec.Range("E1").Value = "Most Recent Update: 11/20/2020`nSee Change Log tab for details"

# Row 36 is the "Intentional self-harm / X83" entry that is being removed
# from the suggested list. Strike through the code/description cells ...
$wsRec.Range("B36:D36").Font.Strikethrough = $true

# ... and append a red explanatory note to the existing notes cell (E36),
# keeping the original black text and adding the new note in red.
$prefix = "8/29 added to address clinicians ability to clearly describe impression; "
$suffix = "11/20/2020 Removed from list, forbidden by the NEMSIS XSD's regex rule."
$e36 = $wsRec.Range("E36")
$e36.Value = $prefix + $suffix
$e36.Characters($prefix.Length + 1, $suffix.Length).Font.Color = 255

# ---------------------------------------------------------------------------
# Change Log sheet
# ---------------------------------------------------------------------------

# Tidy up two existing notes with a trailing period.
$wsLog.Range("F2").Value = "8/21/20 ""NOS"" added to EMS description."
$wsLog.Range("F3").Value = "8/21/20 Removed. R99 already in use for Other: Obvious Death/R99/Ill-defined and unknown cause of mortality."

# Log the new removal as row 4.
$wsLog.Range("A4").Value = 44155
$wsLog.Range("A4").NumberFormat = "m/d/yyyy"
$wsLog.Range("B4").Value = "Emotional State/Behavior"
$wsLog.Range("C4").Value = "Intentional self-harm"
$wsLog.Range("D4").Value = "X83"
$wsLog.Range("E4").Value = "Intentional self-harm by other specified means"
$wsLog.Range("F4").Value = "11/20/2020 Removed from list, forbidden by the NEMSIS XSD's regex rule."

$wsLog.Range("B4:F4").WrapText = $true
$wsLog.Range("A4:F4").VerticalAlignment = -4108
$wsLog.Range("D4:E4").Font.Strikethrough = $true
$wsLog.Rows.Item(4).RowHeight = 30

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping
# ---------------------------------------------------------------------------

$wsLog.Activate() | Out-Null
$wsLog.Range("C9").Select() | Out-Null

$wsRec.Activate() | Out-Null
$wsRec.Range("D29").Select() | Out-Null
